# EncuestasAndalucia2.xlsx - add 4 new poll rows at the top of the data
# table (new Deimos Estadistica, Celeste-Tel, NC Report and Sociometrica
# rows dated 2018-11), and 4 new blank rows at the bottom of the sheet
# (continuing the existing blank-row pattern used for future entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 4 fresh rows above the current row 2 (pushes the existing
#        198 rows of data down to rows 6:202) ------------------------------
$ws.Range("A2:A5").EntireRow.Insert()

# New row 2: Deimos Estadistica, 20-23 Nov 2018
$ws.Range("A2").Value = 43424
$ws.Range("B2").Value = 43427
$ws.Range("C2").Value = 1200
$ws.Range("D2").Value = "Deimos Estadística"
$ws.Range("E2").Value = 18.2
$ws.Range("F2").Value = 36
$ws.Range("G2").Value = 13.2
$ws.Range("H2").Value = 24.9
$ws.Range("I2").Value = 2.8

# New row 3: Celeste-Tel, 19-23 Nov 2018
$ws.Range("A3").Value = 43423
$ws.Range("B3").Value = 43427
$ws.Range("C3").Value = 1000
$ws.Range("D3").Value = "Celeste-Tel"
$ws.Range("E3").Value = 21.2
$ws.Range("F3").Value = 35.8
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 21.6
$ws.Range("I3").Value = 3.2

# New row 4: NC Report, 19-23 Nov 2018
$ws.Range("A4").Value = 43423
$ws.Range("B4").Value = 43427
$ws.Range("C4").Value = 1000
$ws.Range("D4").Value = "NC Report"
$ws.Range("E4").Value = 21.4
$ws.Range("F4").Value = 34.5
$ws.Range("G4").Value = 14.8
$ws.Range("H4").Value = 22.3
$ws.Range("I4").Value = 3.8

# New row 5: Sociometrica, 12-23 Nov 2018
$ws.Range("A5").Value = 43416
$ws.Range("B5").Value = 43427
$ws.Range("C5").Value = 1200
$ws.Range("D5").Value = "Sociométrica"
$ws.Range("E5").Value = 18.6
$ws.Range("F5").Value = 32.8
$ws.Range("G5").Value = 16.2
$ws.Range("H5").Value = 22.7
$ws.Range("I5").Value = 6.5

# --- 2. Append 4 more blank placeholder rows after the old last row -------
#        (old row 198 is now row 202 after the insert above); copy the
#        format of that last blank row down to rows 203:206.
$ws.Range("A202:B202").Copy()
$ws.Range("A203:B206").PasteSpecial(-4122)
$ws.Range("E202:I202").Copy()
$ws.Range("E203:I206").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Match the saved selection / view from the edited workbook ---------
$ws.Range("I4").Select()

Write-Output "done"
